$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New wine entries (Domaine Joseph Colin), appended starting at row 9.
$data = @(
    @(0,    "20Æaine Aligoté Les Jardins de la Cote", "Domaine Joseph Colin", "Bourgogne Générique",    15, 6, 45680),
    @(2021, "Chardonnay Les Hauts de la Combe",        "Domaine Joseph Colin", "Bourgogne Générique",    25, 6, 45680),
    @(2021, "Blanc",                                   "Domaine Joseph Colin", "Chassagne-Montrachet",  47, 6, 45680),
    @(2021, "En Cailleret",                             "Domaine Joseph Colin", "Chassagne-Montrachet",  78, 6, 45680),
    @(2021, "Rouge Vieilles Vignes",                    "Domaine Joseph Colin", "Chassagne-Montrachet",  32, 6, 45680),
    @(2021, "La Garenne",                               "Domaine Joseph Colin", "Puligny-Montrachet",    91, 6, 45680),
    @(2021, "Le Trezin 13",                             "Domaine Joseph Colin", "Puligny-Montrachet",    85, 6, 45680),
    @(2021, "Compendium 135",                           "Domaine Joseph Colin", "Saint-Aubin",           32, 6, 45680),
    @(2021, "Clos du Meix 13",                          "Domaine Joseph Colin", "Saint-Aubin",           43, 6, 45680),
    @(2021, "La Chatenière 135",                        "Domaine Joseph Colin", "Saint-Aubin",           47, 6, 45680)
)

$row = 9
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 7).Value = $entry[4]
    $ws.Cells.Item($row, 10).Value = $entry[5]
    $ws.Cells.Item($row, 11).Value = $entry[6]
    $row = $row + 1
}

# Freeze top row and set the new selection/view state.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B4").Select()
